# Insert a new price record as the first row of the "Tomate" data block
# (row 1072), shifting all subsequent rows down by one. This mirrors the
# source workbook's weekly price-sheet update pattern: a fresh weekly
# observation is prepended and the whole table grows by one row
# (A1:R1182 -> A1:R1183).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 1072:1182 down to 1073:1183, duplicating formatting from the
# row being pushed down (matches Excel's default Insert behaviour).
$ws.Rows("1072:1072").Insert()

# Populate the newly freed row 1072 with the new weekly observation.
$row = 1072
$ws.Cells.Item($row, 1).Value  = 4
$ws.Cells.Item($row, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row, 3).Value  = "Los Lagos"
$ws.Cells.Item($row, 4).Value  = 45166
$ws.Cells.Item($row, 5).Value  = 10
$ws.Cells.Item($row, 6).Value  = 100112020
$ws.Cells.Item($row, 7).Value  = "Tomate"
$ws.Cells.Item($row, 8).Value  = "Larga vida"
$ws.Cells.Item($row, 9).Value  = "Primera"
$ws.Cells.Item($row, 10).Value = 300
$ws.Cells.Item($row, 11).Value = 22000
$ws.Cells.Item($row, 12).Value = 22000
$ws.Cells.Item($row, 13).Value = 22000
$ws.Cells.Item($row, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item($row, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($row, 16).Value = 1222
$ws.Cells.Item($row, 17).Value = 18
$ws.Cells.Item($row, 18).Value = "Hortaliza"
